$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) retains exact text representation (avoid numeric auto-conversion)
$ws.Range("D2:D51").NumberFormat = "@"

# Apply updated cryptocurrency data
$ws.Range("D2").Value = '98.720.24'
$ws.Range("E2").Value = '  -0.20%  '
$ws.Range("D3").Value = '3.345.37'
$ws.Range("E3").Value = '  -1.39%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '259.71'
$ws.Range("E5").Value = '  -0.62%  '
$ws.Range("D6").Value = '647.07'
$ws.Range("E6").Value = '  +2.05%  '
$ws.Range("E7").Value = '  +9.94%  '
$ws.Range("E8").Value = '  +15.71%  '
$ws.Range("D9").Value = '1.10'
$ws.Range("E9").Value = '  +24.40%  '
$ws.Range("E10").Value = '  -0.04%  '
$ws.Range("D11").Value = '3.342.10'
$ws.Range("E11").Value = '  -1.37%  '
$ws.Range("B12").Value = 'Avalanche'
$ws.Range("C12").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D12").Value = '43.99'
$ws.Range("E12").Value = '  +20.81%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '0.208'
$ws.Range("E13").Value = '  +4.05%  '
$ws.Range("E14").Value = '  +7.41%  '
$ws.Range("D15").Value = '100.015.03'
$ws.Range("E15").Value = '  +1.25%  '
$ws.Range("D16").Value = '3.971.49'
$ws.Range("E16").Value = '  -0.85%  '
$ws.Range("D17").Value = '5.57'
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("D18").Value = '3.341.21'
$ws.Range("E18").Value = '  -1.16%  '
$ws.Range("D19").Value = '7.42'
$ws.Range("E19").Value = '  +19.24%  '
$ws.Range("D20").Value = '16.81'
$ws.Range("E20").Value = '  +9.60%  '
$ws.Range("D21").Value = '537.75'
$ws.Range("E21").Value = '  +7.74%  '
$ws.Range("D22").Value = '3.57'
$ws.Range("E22").Value = '  -1.05%  '
$ws.Range("D23").Value = '10.24'
$ws.Range("E23").Value = '  +8.69%  '
$ws.Range("E24").Value = '  -1.04%  '
$ws.Range("E25").Value = '  +52.74%  '
$ws.Range("D26").Value = '103.12'
$ws.Range("E26").Value = '  +13.36%  '
$ws.Range("D27").Value = '6.24'
$ws.Range("E27").Value = '  +7.28%  '
$ws.Range("D28").Value = '12.68'
$ws.Range("E28").Value = '  +4.34%  '
$ws.Range("D29").Value = '3.519.53'
$ws.Range("E29").Value = '  -1.13%  '
$ws.Range("D30").Value = '0.151'
$ws.Range("E30").Value = '  +12.95%  '
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("D32").Value = '10.95'
$ws.Range("E32").Value = '  +13.44%  '
$ws.Range("D33").Value = '0.189'
$ws.Range("E33").Value = '  -7.62%  '
$ws.Range("D34").Value = '0.997'
$ws.Range("E34").Value = '  -0.33%  '
$ws.Range("E35").Value = '  +3.86%  '
$ws.Range("D36").Value = '0.529'
$ws.Range("E36").Value = '  +11.70%  '
$ws.Range("D37").Value = '7.74'
$ws.Range("E37").Value = '  +4.53%  '
$ws.Range("B38").Value = 'PancakeSwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D38").Value = '2.06'
$ws.Range("E38").Value = '  +3.38%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").Value = '0.156'
$ws.Range("E39").Value = '  +2.47%  '
$ws.Range("D40").Value = '515.56'
$ws.Range("E40").Value = '  +2.09%  '
$ws.Range("D41").Value = '24.72'
$ws.Range("E41").Value = '  -0.55%  '
$ws.Range("E42").Value = '  +2.83%  '
$ws.Range("E43").Value = '  +0.60%  '
$ws.Range("D44").Value = '3.33'
$ws.Range("E44").Value = '  -2.00%  '
$ws.Range("D45").Value = '0.811'
$ws.Range("E45").Value = '  +2.30%  '
$ws.Range("D46").Value = '0.0408'
$ws.Range("E46").Value = '  +24.12%  '
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("E48").Value = '  +3.04%  '
$ws.Range("D49").Value = '164.17'
$ws.Range("E49").Value = '  +2.20%  '
$ws.Range("D50").Value = '7.72'
$ws.Range("E50").Value = '  +16.57%  '
$ws.Range("D51").Value = '49.70'
$ws.Range("E51").Value = '  +6.88%  '
